# ANNA-BoM.xlsx edit: "added on off button and remove LCD resstor"
#
# 1. R92 BoM line previously shared its row with an LCD resistor
#    ("R92, LCD-R1"); the LCD resistor reference is removed, leaving just "R92".
# 2. A new BoM entry is added (rows 75-76) for the push/push-push on-off
#    switch (S2 / PUSH SWITCH) including its two reference hyperlinks
#    ("switch" and "cap").
# 3. Minor worksheet cosmetics that came along with the resave: the
#    Description column is widened (and no longer auto bestFit), and the
#    selection/scroll position left at H77 / row 59.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remove "LCD-R1" from the R92 BoM row -------------------------------
$ws.Range("A51").Value = "R92"

# --- 2. Add the new "push switch" BoM row (row 75) + link-only row (76) ----
$ws.Range("A75").Value = "S2"
$ws.Range("B75").Value = "PUSH SWITCH"
$ws.Range("C75").Value = "6pin 2.54mm spacing locked push switch (with white cap)"
$ws.Range("E75").Value = "switch"
$ws.Range("F75").Value = 1

$ws.Range("E76").Value = "cap"
$ws.Range("F76").Value = 1

$ws.Hyperlinks.Add($ws.Range("E75"), "http://www.aliexpress.com/item/6-Pin-DIP-Self-Locking-Push-Button-Switch-Latching-Type-Push-Switch/32714462441.html")
$ws.Hyperlinks.Add($ws.Range("E76"), "http://www.aliexpress.com/item/10pcs-Push-Button-Switch-Cap-White/32714462442.html")

$ws.Range("E75").Style = "Hyperlink"
$ws.Range("E76").Style = "Hyperlink"

# --- 3. Column C width: resave widened it and dropped the bestFit flag -----
$ws.Columns.Item(3).ColumnWidth = 36.6

# --- 4. Leave the selection where the author left it (H77), scrolled down --
$excel.ActiveWindow.ScrollRow = 59
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H77").Select()
